$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Title (paragraph 1, sz=44)
# ------------------------------------------------------------------
$d.Content.Find.Execute("Echoes of Eternity - A Celestial Symphony", $false, $false, $false, $false, $false, $true, 1, $false, "The Profound Insights of Mathematics and Its Widespread Impact", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Author name (paragraph 2, sz=36): "Jaime Torres" -> "Dr. Albert Clayton"
#    (built as 3 logical runs: "Dr", ".", " Albert Clayton")
# ------------------------------------------------------------------
$p2Full = $d.Paragraphs(2).Range
$p2 = $d.Range($p2Full.Start, $p2Full.End - 1)
$p2.Text = "Dr"
$p2.Collapse(0)
$p2.InsertAfter(".")
$p2.Collapse(0)
$p2.InsertAfter(" Albert Clayton")

# ------------------------------------------------------------------
# 3. Email line (paragraph 3, sz=32): jm.torres@umich.edu -> ac.claytonphd@protonmail.com
# ------------------------------------------------------------------
$p3 = $d.Paragraphs(3).Range
$p3.Find.Execute("jm", $false, $false, $false, $false, $false, $true, 1, $false, "ac", 2) | Out-Null
$p3 = $d.Paragraphs(3).Range
$p3.Find.Execute("torres@umich", $false, $false, $false, $false, $false, $true, 1, $false, "claytonphd@protonmail", 2) | Out-Null
$p3 = $d.Paragraphs(3).Range
$p3.Find.Execute("edu", $false, $false, $false, $false, $false, $true, 1, $false, "com", 2) | Out-Null

# ------------------------------------------------------------------
# 4. Body paragraph (paragraph 5, sz=24): full rewrite about Mathematics
# ------------------------------------------------------------------
$para5Items = @(
    @{ t = "text"; v = 'Mathematics, the language of the universe, unravels the intricate patterns and symmetries that permeate existence' },
    @{ t = "text"; v = '.' },
    @{ t = "text"; v = ' It is a boundless realm of exploration, where abstract concepts converge with tangible applications, illuminating the world around us' },
    @{ t = "text"; v = '.' },
    @{ t = "text"; v = ' From the cosmos'' vast expanses to the intricate designs of nature, mathematics provides a framework for comprehending and harnessing the universe''s underlying forces' },
    @{ t = "text"; v = '.' },
    @{ t = "br" },
    @{ t = "br" },
    @{ t = "text"; v = 'Immersed in a world governed by numbers and equations, we discern the rhythmic harmonies of mathematical principles echoing throughout our lives' },
    @{ t = "text"; v = '.' },
    @{ t = "text"; v = ' The Pythagorean theorem unveils the beauty of geometric relationships, guiding architects and engineers in constructing awe-inspiring structures' },
    @{ t = "text"; v = '.' },
    @{ t = "text"; v = ' Calculus, a symphony of change, empowers scientists to model complex phenomena, opening doors to novel technological advancements' },
    @{ t = "text"; v = '.' },
    @{ t = "br" },
    @{ t = "br" },
    @{ t = "text"; v = 'Mathematics serves as a venerable instrument of discovery, propelling humanity''s quest for knowledge' },
    @{ t = "text"; v = '.' },
    @{ t = "text"; v = ' It unveils the secrets of the cosmos, unraveling the mysteries of celestial bodies and guiding astronauts through the vast expanse of space' },
    @{ t = "text"; v = '.' },
    @{ t = "text"; v = ' It unlocks the enigmas of subatomic particles, empowering physicists to explore the fundamental building blocks of matter' },
    @{ t = "text"; v = '.' },
    @{ t = "text"; v = ' Mathematics reveals the intricate machinery of life, enabling biologists to decipher the genetic code and unravel the complexities of the human body' },
    @{ t = "text"; v = '.' }
)

$full5 = $d.Paragraphs(5).Range
$r5 = $d.Range($full5.Start, $full5.End - 1)
$r5.Text = $para5Items[0].v
for ($i = 1; $i -lt $para5Items.Count; $i++) {
    $r5.Collapse(0)
    if ($para5Items[$i].t -eq "br") {
        $r5.InsertAfter([char]11)
    } else {
        $r5.InsertAfter($para5Items[$i].v)
    }
}

# ------------------------------------------------------------------
# 5. Final summary paragraph (last paragraph, default size)
# ------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$finalItems = @(
    'Mathematics, an intellectual odyssey, unveils the universe''s mysteries, propels technological advancements, and underpins our understanding of the cosmos, nature, and life itself',
    '.',
    ' It remains an instrumental force in shaping our world, an enduring testament to the power of human intellect',
    '.'
)
$fullLast = $d.Paragraphs($lastIndex).Range
$rLast = $d.Range($fullLast.Start, $fullLast.End - 1)
$rLast.Text = $finalItems[0]
for ($i = 1; $i -lt $finalItems.Count; $i++) {
    $rLast.Collapse(0)
    $rLast.InsertAfter($finalItems[$i])
}

# ------------------------------------------------------------------
# 6. Add a new empty paragraph at the very end of the document
# ------------------------------------------------------------------
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRange.InsertParagraphAfter()

# ------------------------------------------------------------------
# 7. Fix font typo everywhere: "TimesNewToman" -> "Times New Roman"
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $d.Paragraphs($i).Range.Font.Name = "Times New Roman"
}

Write-Host "Done"
